# Update cryptocurrency price values in column D (Price) on Sheet1
# to reflect the latest scrape as described in the commit message:
# "Updated symbol list on Sat Dec 17 04:41:37 UTC 2022 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "D2"  = "229.99"
    "D3"  = "22.26"
    "D4"  = "5.241"
    "D5"  = "0.05556"
    "D6"  = "3.378"
    "D8"  = "1.064"
    "D9"  = "0.7787"
    "D10" = "0.1370"
    "D11" = "0.07341"
    "D12" = "0.03162"
    "D13" = "0.02947"
    "D14" = "0.09281"
    "D15" = "0.001660"
    "D16" = "3.256"
    "D17" = "0.04781"
    "D18" = "0.0005897"
    "D19" = "0.006224"
    "D20" = "0.005231"
    "D23" = "3.950"
    "D40" = "0.03992"
    "D41" = "0.007135"
    "D42" = "0.1039"
    "D43" = "0.002949"
    "D44" = "0.009748"
    "D45" = "0.00005437"
    "D48" = "0.04446"
}

foreach ($cellRef in $updates.Keys) {
    # Values are stored as text (inline strings) in the workbook, so assign
    # as text to keep the exact formatting (leading/trailing zeros, etc.)
    $ws.Range($cellRef).NumberFormat = "@"
    $ws.Range($cellRef).Value = $updates[$cellRef]
}
